# Update data: 2025-10-29 18:23
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "Metadata": bump the "Last Updated" timestamp
# ---------------------------------------------------------------
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("A2").Value = "29 Oct 2025, 06:23 PM"

# ---------------------------------------------------------------
# Sheet "Top Gainers": refreshed leaderboard rows 42-75
# columns: row, Stock, Latest, Weekly, Monthly
# ---------------------------------------------------------------
$wsGainers = $wb.Worksheets.Item("Top Gainers")
$gainersData = @(
    @(42, "HITECHGEAR", 4.8651, 2.1287, 10.9905),
    @(43, "INDOTHAI", 4.8064, 4.5349, 43.748),
    @(44, "SANDUMA", 4.593, 2.1405, 30.2813),
    @(45, "LLOYDSENT", 4.5646, 1.8339, 11.234),
    @(46, "STAR", 4.5025, 4.4319, 3.662),
    @(47, "RECLTD", 4.4992, 3.4756, 3.4062),
    @(48, "NBCC", 4.4511, 3.1605, 7.6018),
    @(49, "GPPL", 4.4154, 3.4073, 5.0497),
    @(50, "BIL", 4.3654, 9.122199999999999, -0.3203),
    @(51, "HUDCO", 4.3201, 3.8924, 5.3884),
    @(52, "SGMART", 4.2736, 8.258900000000001, 2.5381),
    @(53, "MRPL", 4.2642, 9.7103, 20.0542),
    @(54, "JKIL", 4.1372, 2.9463, 1.7584),
    @(55, "SAMBHV", 4.1349, 2.624, 5.167),
    @(56, "SAPPHIRE", 4.1265, 1.7633, -0.7999000000000001),
    @(57, "PVRINOX", 4.1118, 6.2102, 14.707),
    @(58, "KERNEX", 4.0782, 7.542, 27.2033),
    @(59, "SUNFLAG", 3.997, 4.333, 4.6312),
    @(60, "CMSINFO", 3.9096, 2.6872, 2.8935),
    @(61, "GMBREW", 3.8999, -0.53, 79.029),
    @(62, "GREENLAM", 3.8946, 3.5858, 10.721),
    @(63, "APARINDS", 3.8924, 8.3414, 15.5876),
    @(67, "NPST", 3.7841, -2.0689, -3.5677),
    @(68, "DCW", 3.7544, 2.3219, -3.9753),
    @(69, "RHETAN", 3.754, 4.178, 6.549),
    @(70, "HINDPETRO", 3.6935, 6.9335, 5.7397),
    @(71, "BHARTIHEXA", 3.6718, 7.0877, 15.3332),
    @(72, "HLEGLAS", 3.659, 8.115500000000001, 27.1239),
    @(73, "RHIM", 3.6544, 3.2276, 5.1826),
    @(74, "SHK", 3.6347, 2.388, -1.932),
    @(75, "BCLIND", 3.6271, 2.2945, 0.1728)
)
foreach ($row in $gainersData) {
    $r = $row[0]
    $wsGainers.Cells.Item($r, 2).Value = $row[1]
    $wsGainers.Cells.Item($r, 3).Value = $row[2]
    $wsGainers.Cells.Item($r, 4).Value = $row[3]
    $wsGainers.Cells.Item($r, 5).Value = $row[4]
}

# ---------------------------------------------------------------
# Sheet "1 Month Performance": refreshed leaderboard rows 3-76
# columns: row, Stock, % Change
# ---------------------------------------------------------------
$wsPerf = $wb.Worksheets.Item("1 Month Performance")
$wsPerf.Cells.Item(3, 3).Value = 82.96250000000001

$perfData = @(
    @(5, "PROZONER", 68.3711),
    @(6, "IFBAGRO", 65.2534),
    @(7, "BGRENERGY", 64.929),
    @(8, "ESSARSHPNG", 64.91160000000001),
    @(9, "MAHASTEEL", 56.0982),
    @(10, "INOXGREEN", 51.0181),
    @(11, "STALLION", 46.4325),
    @(12, "ORIENTTECH", 45.3321),
    @(13, "MTARTECH", 40.7213),
    @(14, "TVSSRICHAK", 40.5337),
    @(15, "V2RETAIL", 37.2004),
    @(16, "RAMAPHO", 36.9731),
    @(17, "SANDUMA", 36.9057),
    @(18, "SEJALLTD", 36.8123),
    @(19, "TARACHAND", 36.4813),
    @(20, "NETWEB", 36.1199),
    @(21, "SAMMAANCAP", 35.5128),
    @(22, "ONMOBILE", 35.4702),
    @(23, "SHAREINDIA", 35.3207),
    @(24, "SOUTHBANK", 35.2819),
    @(25, "TVSELECT", 35.1983),
    @(26, "RAMCOSYS", 34.6928),
    @(27, "MAANALU", 34.4803),
    @(28, "MEGASOFT", 33.4399),
    @(29, "BHARATSE", 32.9189),
    @(30, "EMKAY", 30.3743),
    @(31, "ATHERENERG", 29.116),
    @(32, "TATVA", 28.6037),
    @(33, "TERASOFT", 28.3093),
    @(34, "CARTRADE", 27.5713),
    @(35, "ARFIN", 27.4033),
    @(36, "IFBIND", 27.064),
    @(37, "MINDTECK", 26.9415),
    @(38, "BHARATWIRE", 26.5276),
    @(39, "HATSUN", 26.492),
    @(40, "INDORAMA", 26.4516),
    @(41, "ADANIPOWER", 25.8247),
    @(42, "AVALON", 25.7352),
    @(43, "MRPL", 25.6265),
    @(44, "HINDCOPPER", 25.3164),
    @(45, "PRECWIRE", 24.679),
    @(46, "SCI", 24.132),
    @(47, "KICL", 24.1119),
    @(48, "SKYGOLD", 24.1079),
    @(49, "DCBBANK", 23.8922),
    @(50, "AUBANK", 23.6964),
    @(51, "ETHOSLTD", 23.1527),
    @(52, "INDIANB", 22.6689),
    @(53, "PRIVISCL", 22.3984),
    @(54, "CPEDU", 22.3786),
    @(55, "LORDSCHLO", 22.1791),
    @(56, "GUJTHEM", 22.0704),
    @(57, "SURYODAY", 21.8039),
    @(58, "TDPOWERSYS", 21.7743),
    @(59, "ORBTEXP", 21.6115),
    @(60, "CEATLTD", 20.0239),
    @(61, "ATL", 19.9362),
    @(62, "HITECHGEAR", 19.8096),
    @(63, "GRMOVER", 19.7859),
    @(64, "FEDERALBNK", 19.6872),
    @(65, "SUBROS", 19.6508),
    @(66, "USHAMART", 19.6172),
    @(67, "BANKINDIA", 19.3067),
    @(68, "RBLBANK", 19.2556),
    @(69, "MOLDTECH", 19.1891),
    @(70, "THOMASCOTT", 19.1649),
    @(71, "KARURVYSYA", 19.11),
    @(72, "IIFL", 18.9853),
    @(73, "LUMAXIND", 18.8057),
    @(74, "REPRO", 18.689),
    @(75, "TINNARUBR", 18.5804),
    @(76, "SRM", 18.4636)
)
foreach ($row in $perfData) {
    $r = $row[0]
    $wsPerf.Cells.Item($r, 2).Value = $row[1]
    $wsPerf.Cells.Item($r, 3).Value = $row[2]
}

Write-Host "Edit complete"
